# New Submission Synced: 2026-02-09 17:43:43
# Appends the new Google-Form submission row to the "JSS 3B" response sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("JSS 3B")

$ws.Range("A8").Value = "2026-02-09 17:43:43"
$ws.Range("B8").Value = "Ishaku Yusuf Dawha "

# Admission No is stored as text (it is sometimes non-numeric, e.g. "Number 3"),
# so force text formatting before writing the numeric-looking "36", then drop
# back to the sheet's normal (unformatted) style so only the cell's stored
# type - not its appearance - changes.
$ws.Range("C8").NumberFormat = "@"
$ws.Range("C8").Value = "36"
$ws.Range("C8").Style = "Normal"

$ws.Range("D8").Value = 7
